$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Insertion 1: new "Every entry..." / "Solution: " / blank
# paragraphs, inserted right after the existing blank paragraph
# that separates the Batch Processor section from the Waveforms
# Server header.
# ------------------------------------------------------------------
$pBlank1 = $d.Paragraphs.Item(10)
$pBlank1.Range.InsertParagraphAfter()
$pBlank1.Range.InsertParagraphAfter()
$pBlank1.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(11)
$pNew2 = $d.Paragraphs.Item(12)

$pNew1.Range.Text = "Every entry in the database is being reprocessed by the batch processor."
$pNew1.Range.Font.Color = 255

$pNew2.Range.Text = "Solution: "

# ------------------------------------------------------------------
# Insertion 2: new "I need to access/modify the network
# configuration." / "Solution: Open the UniFi control program..."
# paragraphs, inserted right after the "Ubiquiti Unifi:" header.
# ------------------------------------------------------------------
$pUbiquitiHeader = $d.Paragraphs.Item(19)
$pUbiquitiHeader.Range.InsertParagraphAfter()
$pUbiquitiHeader.Range.InsertParagraphAfter()
$pNew3 = $d.Paragraphs.Item(20)
$pNew4 = $d.Paragraphs.Item(21)

$pNew3.Range.Text = "I need to access/modify the network configuration."
$pNew3.Range.Font.Color = 255

$pNew4.Range.Text = "Solution: Open the UniFi control program. There should be a shortcut on the desktop. (This shortcut points to the rather non-obvious " + [char]34 + "C:\Users\dusty\Ubiquiti UniFi\lib\ace.jar" + [char]34 + " ui). It" + [char]8217 + "ll take a while to connect, but once it does, select " + [char]8220 + "launch a browser to manage the network." + [char]8221 + " This must be done from dusty, although other machines can be used if you " + [char]8216 + "adopt" + [char]8217 + " the UniFi controller from them. You should avoid doing that unless necessary, as it" + [char]8217 + "s a bit of a process."
